$wb = $excel.ActiveWorkbook

# Update the "59e52cfe-..." row datetimes on the zh-cn sheet (rows 3 and 5
# both reference the same shared-string values, so both must be updated)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-21 04:21:49"
$wsZhCn.Range("H3").Value = "2016-03-21 04:22:12"
$wsZhCn.Range("E5").Value = "2016-03-21 04:21:49"
$wsZhCn.Range("H5").Value = "2016-03-21 04:22:12"

# Update the "59e52cfe-..." row datetimes on the de-de sheet (rows 3 and 5
# both reference the same shared-string values, so both must be updated)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-21 04:21:53"
$wsDeDe.Range("H3").Value = "2016-03-21 04:22:17"
$wsDeDe.Range("E5").Value = "2016-03-21 04:21:53"
$wsDeDe.Range("H5").Value = "2016-03-21 04:22:17"
